# Apply "Camera controls added but models buggy" edits to the rubric sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11 (10. Directional light w/ specular highlights): EARNED score drops
# to 0, with a note in column D explaining why.
$ws.Range("C11").Value = 0
$ws.Range("D11").Value = "Note: specular lights must react to camera movements"

# Rows that are "in progress" get an "I want to do this" marker in column D.
$ws.Range("D18").Value = "I want to do this"
$ws.Range("D19").Value = "I want to do this"
$ws.Range("D20").Value = "I want to do this"
$ws.Range("D33").Value = "I want to do this"
$ws.Range("D51").Value = "I want to do this"
$ws.Range("D53").Value = "I want to do this"

# Row 52 (commit regularly) is actively being worked on.
$ws.Range("D52").Value = "Doing this"

# Move the active selection to D8, matching the edited workbook view.
$ws.Range("D8").Select()
